$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# --- Text / label updates (row 4: XLK/Technology -> XLI/Industrial) ---
$ws.Range("A4").Value = "XLI"
$ws.Range("B4").Value = "Industrial Select Sector SPDR Fund"
$ws.Range("C4").Value = "Industrials Funds"

# --- Update "as of" date in confidential disclosure text (row 9) ---
$ws.Range("A9").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-21 for illustrative purposes only and are subject to change."

# --- Update numeric Weight (D) / Percent Change (E) values ---
$ws.Range("D2").Value = 0.2517476653662153
$ws.Range("E2").Value = 0.01816173763651996

$ws.Range("D3").Value = 0.2469753458585881
$ws.Range("E3").Value = 0.01385681293302521

$ws.Range("D4").Value = 0.2467070937662453
$ws.Range("E4").Value = 0.01365324766589704

$ws.Range("D5").Value = 0.2545698950089512
$ws.Range("E5").Value = 0.001710526315789496

$ws.Range("D6").Value = 0.9999999999999999
$ws.Range("E6").Value = 0.01179826777236381

# Restore sheet protection (original password hash cannot be reconstructed
# from the legacy hash stored in the source file, so re-protect with the
# same effective protected state the sheet had before editing).
$ws.Protect()
